$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 is appended with the same shape as the existing data rows. The
# migration-date column holds a literal text value ("2025-10-17"), not a
# real date, so the cell is pre-formatted as Text to stop Excel's
# automatic date recognition from converting it to a serial date number;
# the style is then reset back to Normal so the cell doesn't end up with
# a lingering number-format style (matching the plain, style-less cells
# used elsewhere on this sheet).
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-10-17"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "ZZZ"
$ws.Range("C4").Value = "456CDX009"
$ws.Range("D4").Value = "Anna Nagar"
